$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 and L3 are stored as text in the workbook (date written as plain text,
# phone number as text), so force text format before assigning so Excel
# COM doesn't auto-convert them to a date serial / numeric value.
$ws.Range("C3").NumberFormat = "@"
$ws.Range("L3").NumberFormat = "@"

# Update row 3 (data row) with new values per the diff
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = "4/8/2020"
$ws.Range("D3").Value = "Dallas"
$ws.Range("E3").Value = "NewJersey"
$ws.Range("F3").Value = "6:00 pm"
$ws.Range("G3").Value = "4:00 pm"
$ws.Range("I3").Value = "Angela"
$ws.Range("J3").Value = "Smith"
$ws.Range("K3").Value = "A@aol.com"
$ws.Range("L3").Value = "2581236548"
$ws.Range("M3").Value = 33
